$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "statut" column (A) and "statut_label" column (B) for row 2
# were changed from the black-square/"noir" entry to the blue-book/"bleu" entry.
$ws.Range("A2").Value = "📘"
$ws.Range("B2").Value = "bleu"
